# Apply the edit described in the diff:
# - Insert two new rows at position 151 (this shifts the old rows 151..189
#   down to 153..191, which already reproduces every value change the diff
#   shows for those shifted rows).
# - Populate the two newly inserted rows (151 and 152) with their data.
# - The worksheet dimension will automatically grow from A1:T189 to A1:T191.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows before the current row 151; everything below shifts down by 2.
$ws.Rows("151:152").Insert()

# --- New row 151 ---
$ws.Range("A151").Value = 5
$ws.Range("B151").Value = "Macroferia Regional de Talca"
$ws.Range("C151").Value = "Maule"
$ws.Range("D151").Value2 = 44985
$ws.Range("E151").Value = 7
$ws.Range("F151").Value = "Fruta"
$ws.Range("G151").Value = 100103
$ws.Range("H151").Value = "Frutos de hueso (carozo)"
$ws.Range("I151").Value = 100103002
$ws.Range("J151").Value = "Ciruela"
$ws.Range("K151").Value = "Black Amber"
$ws.Range("L151").Value = "Primera"
$ws.Range("M151").Value = 250
$ws.Range("N151").Value = 8000
$ws.Range("O151").Value = 8000
$ws.Range("P151").Value = 8000
$ws.Range("Q151").Value = '$/bandeja 18 kilos granel'
$ws.Range("R151").Value = "Provincia de Curicó"
$ws.Range("S151").Value = 444
$ws.Range("T151").Value = 18

# --- New row 152 ---
$ws.Range("A152").Value = 5
$ws.Range("B152").Value = "Macroferia Regional de Talca"
$ws.Range("C152").Value = "Maule"
$ws.Range("D152").Value2 = 44985
$ws.Range("E152").Value = 7
$ws.Range("F152").Value = "Fruta"
$ws.Range("G152").Value = 100103
$ws.Range("H152").Value = "Frutos de hueso (carozo)"
$ws.Range("I152").Value = 100103002
$ws.Range("J152").Value = "Ciruela"
$ws.Range("K152").Value = "Black Amber"
$ws.Range("L152").Value = "Segunda"
$ws.Range("M152").Value = 180
$ws.Range("N152").Value = 5000
$ws.Range("O152").Value = 5000
$ws.Range("P152").Value = 5000
$ws.Range("Q152").Value = '$/bandeja 18 kilos granel'
$ws.Range("R152").Value = "Provincia de Curicó"
$ws.Range("S152").Value = 278
$ws.Range("T152").Value = 18
